$d = $word.ActiveDocument

# --- 1. Fix last-name ordering: "Armenta Uribe Javier" -> "Uribe Armenta Javier" ---
# First do a straightforward text replace (this temporarily lands as one run).
$findRange = $d.Content
$findRange.Find.Execute("Armenta Uribe Javier", $false, $false, $false, $false, $false, $true, 1, $false, "Uribe Armenta Javier", 2) | Out-Null

# Re-locate the corrected text so we can split it back into three runs, same as
# Word naturally does when a sentence gets edited word-by-word: "Uribe ", then
# "Armenta ", then "Javier" - each its own run, all sharing identical rPr.
$target = $d.Content
$target.Find.Execute("Uribe Armenta Javier") | Out-Null
$start = $target.Start
$uribeEnd = $start + 6      # "Uribe " is 6 characters
$armentaEnd = $uribeEnd + 8 # "Armenta " is 8 characters

$r1 = $d.Range($start, $uribeEnd)
$r2 = $d.Range($uribeEnd, $armentaEnd)

# Toggling Bold on then off forces a run split at each boundary without
# leaving any actual formatting difference behind in the saved XML.
$r1.Font.Bold = 1
$r1.Font.Bold = 0
$r2.Font.Bold = 1
$r2.Font.Bold = 0

# --- 2. Mark the two picture-containing runs as NoProof ---
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $p.Range.NoProofing = 1
    }
}
